# "Range Migration vs Resolution.xlsx" - HW2
# Commit: code run-time speed improvement
#   1. Use fft with zero padding to speed up calculation
#   2. pre-allocate memory to speed up (most of the improvement)
#
# The actual spreadsheet change (the code itself lives outside the workbook):
# insert a small header block above the existing results table that labels
# columns E/F ("Range Migration" / "Range Resolution"), pushing the original
# header/data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row at row 4 (shifts old row4 header -> row5, old row5 data -> row6,
# old row7 -> row8). Row 3 stays blank except for a spacer cell below B2.
$ws.Range("A4:J4").EntireRow.Insert()

# Give the (now empty) B3 spacer cell the same numeric style as B2 (0.00, centered)
# so it keeps the style slot already present in the workbook instead of creating a
# new one.
$ws.Range("B3").NumberFormat = "0.00"
$ws.Range("B3").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# New labels sitting above the "Difference" calc columns (E/F) in the results table.
$ws.Range("E4").Value = "Range Migration"
$ws.Range("F4").Value = "Range Resolution"

# Auto-fit the two columns so the new header text is fully visible (this is what
# flips the columns from the old shared 8.55-wide pair to their own best-fit widths).
$ws.Columns("E:F").EntireColumn.AutoFit()

# Match the author's final selection in the saved file.
$ws.Range("G6").Select()
